$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "week"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "provider"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "state"

# ---- week sheet: weekly Sunday dates, formatted as short dates ----
$dates = @(43975,43982,43989,43996,44003,44010,44017,44024,44031,44038,44045,44052,44059,44066,44073,44080,44087,44094,44101,44108,44115,44122,44129,44136,44143,44150,44157,44164,44171,44178,44185,44192)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $dates[$i]
}

# Apply the date number format to A1 first, then propagate the resulting
# style (not the value!) to the rest of the column so every cell shares one xf.
$ws1.Cells.Item(1, 1).NumberFormat = "mm-dd-yy"
$ws1.Cells.Item(1, 1).Copy()
$ws1.Range("A2:A32").PasteSpecial(-4122)  # xlPasteFormats

# ---- provider sheet: provider codes, kept as zero-padded text ----
$providers = @("015009","015010","015012","015014","015015","015016","015019","015023","015024","015027","015028","015031","015032","015034","015035","015037","015040")

for ($i = 0; $i -lt $providers.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Formula = '=TEXT(' + [int]$providers[$i] + ',"000000")'
}
$provRange = $ws2.Range("A1:A" + $providers.Length)
$provRange.Copy()
$provRange.PasteSpecial(-4163)

# ---- state sheet: single state code ----
$ws3.Cells.Item(1, 1).Formula = '=TEXT("AL","@")'
$stateRange = $ws3.Range("A1:A1")
$stateRange.Copy()
$stateRange.PasteSpecial(-4163)

# restore selection/active sheet to match original workbook view
$ws1.Range("A1").Select()
$ws1.Activate()
